$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "24.889.09"
$rng.ClearFormats()
$rng = $ws.Range("E2")
$rng.NumberFormat = "@"
$rng.Value = "  +2.33%  "
$rng.ClearFormats()
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "1.707.00"
$rng.ClearFormats()
$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = "  +1.93%  "
$rng.ClearFormats()
$rng = $ws.Range("D4")
$rng.NumberFormat = "@"
$rng.Value = "1.000"
$rng.ClearFormats()
$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$rng.Value = "  -0.19%  "
$rng.ClearFormats()
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "315.46"
$rng.ClearFormats()
$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$rng.Value = "  +0.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = "0.9999"
$rng.ClearFormats()
$rng = $ws.Range("E6")
$rng.NumberFormat = "@"
$rng.Value = "  -0.25%  "
$rng.ClearFormats()
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "0.3949"
$rng.ClearFormats()
$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$rng.Value = "  +1.72%  "
$rng.ClearFormats()
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = "0.4017"
$rng.ClearFormats()
$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$rng.Value = "  +0.85%  "
$rng.ClearFormats()
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "1.484"
$rng.ClearFormats()
$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = "  +1.46%  "
$rng.ClearFormats()
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "52.86"
$rng.ClearFormats()
$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = "  +1.15%  "
$rng.ClearFormats()
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "1.001"
$rng.ClearFormats()
$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$rng.Value = "  -0.13%  "
$rng.ClearFormats()
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "0.08826"
$rng.ClearFormats()
$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$rng.Value = "  +1.44%  "
$rng.ClearFormats()
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "25.93"
$rng.ClearFormats()
$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$rng.Value = "  +4.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "7.473"
$rng.ClearFormats()
$rng = $ws.Range("E14")
$rng.NumberFormat = "@"
$rng.Value = "  +0.61%  "
$rng.ClearFormats()
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "7.990"
$rng.ClearFormats()
$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$rng.Value = "  +0.94%  "
$rng.ClearFormats()
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "0.00001348"
$rng.ClearFormats()
$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$rng.Value = "  +1.01%  "
$rng.ClearFormats()
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "1.765.39"
$rng.ClearFormats()
$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$rng.Value = "  +5.63%  "
$rng.ClearFormats()
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "96.49"
$rng.ClearFormats()
$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$rng.Value = "  -1.31%  "
$rng.ClearFormats()
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "0.07185"
$rng.ClearFormats()
$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = "  +1.67%  "
$rng.ClearFormats()
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "20.59"
$rng.ClearFormats()
$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$rng.Value = "  +5.57%  "
$rng.ClearFormats()
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "7.379"
$rng.ClearFormats()
$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = "  +2.43%  "
$rng.ClearFormats()
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "0.9992"
$rng.ClearFormats()
$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$rng.Value = "  -0.30%  "
$rng.ClearFormats()
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "14.45"
$rng.ClearFormats()
$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = "  +2.83%  "
$rng.ClearFormats()
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "24.884.46"
$rng.ClearFormats()
$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$rng.Value = "  +2.44%  "
$rng.ClearFormats()
$rng = $ws.Range("B25")
$rng.NumberFormat = "@"
$rng.Value = "LidoDAOToken"
$rng.ClearFormats()
$rng = $ws.Range("C25")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$rng.ClearFormats()
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "2.985"
$rng.ClearFormats()
$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$rng.Value = "  +0.55%  "
$rng.ClearFormats()
$rng = $ws.Range("B26")
$rng.NumberFormat = "@"
$rng.Value = "Toncoin"
$rng.ClearFormats()
$rng = $ws.Range("C26")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$rng.ClearFormats()
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "2.358"
$rng.ClearFormats()
$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$rng.Value = "  +1.16%  "
$rng.ClearFormats()
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "23.76"
$rng.ClearFormats()
$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$rng.Value = "  +6.31%  "
$rng.ClearFormats()
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "6.231"
$rng.ClearFormats()
$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = "  +19.54%  "
$rng.ClearFormats()
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "160.89"
$rng.ClearFormats()
$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = "  -2.08%  "
$rng.ClearFormats()
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "149.46"
$rng.ClearFormats()
$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$rng.Value = "  +9.33%  "
$rng.ClearFormats()
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "8.498"
$rng.ClearFormats()
$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$rng.Value = "  -1.98%  "
$rng.ClearFormats()
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "2.453"
$rng.ClearFormats()
$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$rng.Value = "  +24.26%  "
$rng.ClearFormats()
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "1.939.41"
$rng.ClearFormats()
$rng = $ws.Range("E33")
$rng.NumberFormat = "@"
$rng.Value = "  +4.38%  "
$rng.ClearFormats()
$rng = $ws.Range("B34")
$rng.NumberFormat = "@"
$rng.Value = "InternetComputer(DFINITY)"
$rng.ClearFormats()
$rng = $ws.Range("C34")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$rng.ClearFormats()
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "7.323"
$rng.ClearFormats()
$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$rng.Value = "  +0.52%  "
$rng.ClearFormats()
$rng = $ws.Range("B35")
$rng.NumberFormat = "@"
$rng.Value = "Hedera"
$rng.ClearFormats()
$rng = $ws.Range("C35")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$rng.ClearFormats()
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "0.08571"
$rng.ClearFormats()
$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$rng.Value = "  -1.51%  "
$rng.ClearFormats()
$rng = $ws.Range("B36")
$rng.NumberFormat = "@"
$rng.Value = "ImmutableX"
$rng.ClearFormats()
$rng = $ws.Range("C36")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$rng.ClearFormats()
$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = "1.044"
$rng.ClearFormats()
$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$rng.Value = "  +1.41%  "
$rng.ClearFormats()
$rng = $ws.Range("B37")
$rng.NumberFormat = "@"
$rng.Value = "VeChain"
$rng.ClearFormats()
$rng = $ws.Range("C37")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$rng.ClearFormats()
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "0.03120"
$rng.ClearFormats()
$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$rng.Value = "  +7.23%  "
$rng.ClearFormats()
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "0.2849"
$rng.ClearFormats()
$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$rng.Value = "  +3.43%  "
$rng.ClearFormats()
$rng = $ws.Range("B39")
$rng.NumberFormat = "@"
$rng.Value = "FraxShare"
$rng.ClearFormats()
$rng = $ws.Range("C39")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$rng.ClearFormats()
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = "10.84"
$rng.ClearFormats()
$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$rng.Value = "  +2.06%  "
$rng.ClearFormats()
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "0.09485"
$rng.ClearFormats()
$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$rng.Value = "  +4.62%  "
$rng.ClearFormats()
$rng = $ws.Range("B41")
$rng.NumberFormat = "@"
$rng.Value = "TheSandbox"
$rng.ClearFormats()
$rng = $ws.Range("C41")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$rng.ClearFormats()
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "0.8302"
$rng.ClearFormats()
$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$rng.Value = "  +6.38%  "
$rng.ClearFormats()
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "14.07"
$rng.ClearFormats()
$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = "  +0.55%  "
$rng.ClearFormats()
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "1.479"
$rng.ClearFormats()
$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$rng.Value = "  +1.42%  "
$rng.ClearFormats()
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = "17.49"
$rng.ClearFormats()
$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$rng.Value = "  +3.46%  "
$rng.ClearFormats()
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "2.673"
$rng.ClearFormats()
$rng = $ws.Range("E45")
$rng.NumberFormat = "@"
$rng.Value = "  +4.36%  "
$rng.ClearFormats()
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "0.7398"
$rng.ClearFormats()
$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$rng.Value = "  +4.01%  "
$rng.ClearFormats()
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "4.251"
$rng.ClearFormats()
$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = "  +1.15%  "
$rng.ClearFormats()
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "1.383"
$rng.ClearFormats()
$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$rng.Value = "  -1.34%  "
$rng.ClearFormats()
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "0.08667"
$rng.ClearFormats()
$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$rng.Value = "  +8.54%  "
$rng.ClearFormats()
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = "0.9999"
$rng.ClearFormats()
$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$rng.Value = "  -0.24%  "
$rng.ClearFormats()
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "139.48"
$rng.ClearFormats()
$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$rng.Value = "  +0.86%  "
$rng.ClearFormats()

Write-Host "Done applying 116 cell updates."
